$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5234368443489075
$ws.Range("B1").Value = 2.185462951660156
$ws.Range("C1").Value = 7.777897834777832
$ws.Range("D1").Value = 2.184454441070557
$ws.Range("E1").Value = 0.9811344742774963
